# Disease incidence.xlsx edit: add an "AVE" (average) summary row and a
# MIN() check cell to the "DI" sheet, then restore the selections that were
# left active on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet 1")
$ws2 = $wb.Worksheets.Item("DI")

# --- DI sheet: new row 7 ("AVE" site averages) and row 8 (overall MIN) ---
$ws2.Range("A7").Value = "AVE"
$ws2.Range("B7").Formula = "=AVERAGE(B`$2:B`$6)"
$ws2.Range("C7:F7").Formula = "=AVERAGE(C`$2:C`$6)"
$ws2.Range("B8").Formula = "=MIN(B2:F6)"

# --- Column widths on DI (best-fit recalculated once column A/D appear) ---
$ws2.Columns.Item(1).ColumnWidth = 3.5
$ws2.Columns.Item(2).ColumnWidth = 13.5
$ws2.Columns.Item(3).ColumnWidth = 7.8333333
$ws2.Columns.Item(4).ColumnWidth = 7.6666667
$ws2.Columns.Item(5).ColumnWidth = 11.6666667
$ws2.Columns.Item(6).ColumnWidth = 7.5

# --- Restore the cursor/selection state recorded in the saved file ---
$ws1.Range("F22").Select() | Out-Null
$ws2.Range("J5").Select() | Out-Null
